$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.766253692609723
$ws.Range("C2").Value = 0.03487257306364455
$ws.Range("D2").Value = 0.4923411594385243
$ws.Range("E2").Value = 0.1466535952816503
$ws.Range("G2").Value = 3.211953454185817
$ws.Range("H2").Value = 2.483316680697243
$ws.Range("J2").Value = 0.0487080606150796
$ws.Range("K2").Value = 1.129349248461011
$ws.Range("L2").Value = 0.4266135680399188
$ws.Range("M2").Value = 0.4551536631419779
$ws.Range("N2").Value = 4.266766516072522
$ws.Range("B3").Value = 1.736537289183218
$ws.Range("C3").Value = 0.03203795572482449
$ws.Range("D3").Value = 0.4911562912160434
$ws.Range("E3").Value = 0.1470383042868111
$ws.Range("G3").Value = 3.209052450671763
$ws.Range("H3").Value = 2.487718045696269
$ws.Range("J3").Value = 0.0485000923745762
$ws.Range("K3").Value = 1.098887521905539
$ws.Range("L3").Value = 0.4255283808693662
$ws.Range("M3").Value = 0.4502797440044226
$ws.Range("N3").Value = 4.284718838245013
$ws.Range("B4").Value = 1.719194399606693
$ws.Range("C4").Value = 0.03028105029029859
$ws.Range("D4").Value = 0.4906163870430191
$ws.Range("E4").Value = 0.1473111995798444
$ws.Range("G4").Value = 3.208465141679397
$ws.Range("H4").Value = 2.491188917446152
$ws.Range("J4").Value = 0.04836854933486467
$ws.Range("K4").Value = 1.080813959041166
$ws.Range("L4").Value = 0.4250258737939703
$ws.Range("M4").Value = 0.4475014854760353
$ws.Range("N4").Value = 4.296641801196074
$ws.Range("B5").Value = 1.712354527909696
$ws.Range("C5").Value = 0.02956093543991045
$ws.Range("D5").Value = 0.4904436309721802
$ws.Range("E5").Value = 0.1474316515133456
$ws.Range("G5").Value = 3.208526125980683
$ws.Range("H5").Value = 2.492796740605911
$ws.Range("J5").Value = 0.04831397236675805
$ws.Range("K5").Value = 1.073607499737847
$ws.Range("L5").Value = 0.4248623658032074
$ws.Range("M5").Value = 0.4464233202725687
$ws.Range("N5").Value = 4.301726925897228
$ws.Range("B6").Value = 1.711232524675523
$ws.Range("C6").Value = 0.02944110875230166
$ws.Range("D6").Value = 0.4904178017665828
$ws.Range("E6").Value = 0.1474522114468222
$ws.Range("G6").Value = 3.208554393780219
$ws.Range("H6").Value = 2.493075405841992
$ws.Range("J6").Value = 0.04830485107019999
$ws.Range("K6").Value = 1.072420462489703
$ws.Range("L6").Value = 0.4248377102123797
$ws.Range("M6").Value = 0.4462475560063019
$ws.Range("N6").Value = 4.302584984247957
$ws.Range("B7").Value = 1.719101233042522
$ws.Range("C7").Value = 0.03027135545625015
$ws.Range("D7").Value = 0.4906138657310493
$ws.Range("E7").Value = 0.1473127865766788
$ws.Range("G7").Value = 3.208464748039077
$ws.Range("H7").Value = 2.491209817738095
$ws.Range("J7").Value = 0.04836781723229588
$ws.Range("K7").Value = 1.080716127477729
$ws.Range("L7").Value = 0.4250235014722819
$ws.Range("M7").Value = 0.4474867261976279
$ws.Range("N7").Value = 4.296709464123737
$ws.Range("B8").Value = 1.755820265301594
$ws.Range("C8").Value = 0.03389859136862583
$ws.Range("D8").Value = 0.4918937268937071
$ws.Range("E8").Value = 0.146778641869977
$ws.Range("G8").Value = 3.210705439996559
$ws.Range("H8").Value = 2.484674897224835
$ws.Range("J8").Value = 0.04863714968727884
$ws.Range("K8").Value = 1.118715400878415
$ws.Range("L8").Value = 0.4262054439387697
$ws.Range("M8").Value = 0.4534287040597818
$ws.Range("N8").Value = 4.272769687159752
$ws.Range("B9").Value = 1.83498017124586
$ws.Range("C9").Value = 0.04088270811575967
$ws.Range("D9").Value = 0.4958890905497242
$ws.Range("E9").Value = 0.1460213732962998
$ws.Range("G9").Value = 3.224571901202864
$ws.Range("H9").Value = 2.477950054636011
$ws.Range("J9").Value = 0.04913496063235279
$ws.Range("K9").Value = 1.198226448671107
$ws.Range("L9").Value = 0.4298200363279392
$ws.Range("M9").Value = 0.4667784962644674
$ws.Range("N9").Value = 4.232963499863871
$ws.Range("B10").Value = 1.89749341975579
$ws.Range("C10").Value = 0.0459379430981528
$ws.Range("D10").Value = 0.4997269535951858
$ws.Range("E10").Value = 0.1456408065826906
$ws.Range("G10").Value = 3.240538493448753
$ws.Range("H10").Value = 2.476714854509055
$ws.Range("J10").Value = 0.04948249076918465
$ws.Range("K10").Value = 1.2596898890516
$ws.Range("L10").Value = 0.4332631748022493
$ws.Range("M10").Value = 0.4776189528392365
$ws.Range("N10").Value = 4.208066973904096
$ws.Range("B11").Value = 1.92687705062508
$ws.Range("C11").Value = 0.04822177658047622
$ws.Range("D11").Value = 0.5016682600411855
$ws.Range("E11").Value = 0.1455056218303667
$ws.Range("G11").Value = 3.249058284178489
$ws.Range("H11").Value = 2.476956242331227
$ws.Range("J11").Value = 0.0496366998002209
$ws.Range("K11").Value = 1.288313720702064
$ws.Range("L11").Value = 0.4349999594147249
$ws.Range("M11").Value = 0.4827743153133426
$ws.Range("N11").Value = 4.197684652569492
$ws.Range("B12").Value = 1.938139658951968
$ws.Range("C12").Value = 0.04908436690975293
$ws.Range("D12").Value = 0.5024314206779223
$ws.Range("E12").Value = 0.1454598667635381
$ws.Range("G12").Value = 3.252465215845774
$ws.Range("H12").Value = 2.477163024924863
$ws.Range("J12").Value = 0.04969454075870949
$ws.Range("K12").Value = 1.299248186429651
$ws.Range("L12").Value = 0.4356820869569589
$ws.Range("M12").Value = 0.4847586579141137
$ws.Range("N12").Value = 4.193888745937215
$ws.Range("B13").Value = 1.935708026912891
$ws.Range("C13").Value = 0.04889869219468324
$ws.Range("D13").Value = 0.5022658149167825
$ws.Range("E13").Value = 0.1454694794457332
$ws.Range("G13").Value = 3.25172343687359
$ws.Range("H13").Value = 2.477113361443742
$ws.Range("J13").Value = 0.04968210828001851
$ws.Range("K13").Value = 1.296889021029727
$ws.Range("L13").Value = 0.4355340925064866
$ws.Range("M13").Value = 0.4843298671640639
$ws.Range("N13").Value = 4.194700228811726
$ws.Range("B14").Value = 1.927800915390549
$ws.Range("C14").Value = 0.0482927873018042
$ws.Range("D14").Value = 0.5017304843864991
$ws.Range("E14").Value = 0.1455017486858203
$ws.Range("G14").Value = 3.249334953620746
$ws.Range("H14").Value = 2.476970942963987
$ws.Range("J14").Value = 0.04964146950455728
$ws.Range("K14").Value = 1.289211398193771
$ws.Range("L14").Value = 0.4350555889467245
$ws.Range("M14").Value = 0.482936925178258
$ws.Range("N14").Value = 4.197369642241625
$ws.Range("B15").Value = 1.922975237820594
$ws.Range("C15").Value = 0.04792136112070011
$ws.Range("D15").Value = 0.501406227116135
$ws.Range("E15").Value = 0.1455222219593253
$ws.Range("G15").Value = 3.247895467172441
$ws.Range("H15").Value = 2.476898728643533
$ws.Range("J15").Value = 0.04961650496377779
$ws.Range("K15").Value = 1.284521029579508
$ws.Range("L15").Value = 0.4347656730692506
$ws.Range("M15").Value = 0.4820878881432833
$ws.Range("N15").Value = 4.199022401674171
$ws.Range("B16").Value = 1.895592135023094
$ws.Range("C16").Value = 0.045788373194938
$ws.Range("D16").Value = 0.4996040110888345
$ws.Range("E16").Value = 0.1456504030003849
$ws.Range("G16").Value = 3.240006988067762
$ws.Range("H16").Value = 2.476715231611564
$ws.Range("J16").Value = 0.04947233501933113
$ws.Range("K16").Value = 1.257832593550035
$ws.Range("L16").Value = 0.4331530963724788
$ws.Range("M16").Value = 0.4772865372033124
$ws.Range("N16").Value = 4.208764464665521
$ws.Range("B17").Value = 1.879035524308165
$ws.Range("C17").Value = 0.04447582180483778
$ws.Range("D17").Value = 0.4985484144233823
$ws.Range("E17").Value = 0.1457387410410576
$ws.Range("G17").Value = 3.235489473152597
$ws.Range("H17").Value = 2.476808304358997
$ws.Range("J17").Value = 0.04938289936495366
$ws.Range("K17").Value = 1.241629966582394
$ws.Range("L17").Value = 0.4322074474083877
$ws.Range("M17").Value = 0.474398367375457
$ws.Range("N17").Value = 4.214982500393759
$ws.Range("B18").Value = 1.869601667661641
$ws.Range("C18").Value = 0.04371938646838203
$ws.Range("D18").Value = 0.4979596636013781
$ws.Range("E18").Value = 0.1457931226675573
$ws.Range("G18").Value = 3.233009386952659
$ws.Range("H18").Value = 2.476937454453292
$ws.Range("J18").Value = 0.04933109251419587
$ws.Range("K18").Value = 1.232373136990958
$ws.Range("L18").Value = 0.4316795860309099
$ws.Range("M18").Value = 0.4727582562687402
$ws.Range("N18").Value = 4.218647730288438
$ws.Range("B19").Value = 1.866422834881973
$ws.Range("C19").Value = 0.04346301365016814
$ws.Range("D19").Value = 0.4977634852659065
$ws.Range("E19").Value = 0.1458121494089273
$ws.Range("G19").Value = 3.23218998564181
$ws.Range("H19").Value = 2.476994174017221
$ws.Range("J19").Value = 0.04931348865524576
$ws.Range("K19").Value = 1.229249673286006
$ws.Range("L19").Value = 0.4315036201528386
$ws.Range("M19").Value = 0.4722065674348741
$ws.Range("N19").Value = 4.219903961800952
$ws.Range("B20").Value = 1.880788788311634
$ws.Range("C20").Value = 0.04461569920799491
$ws.Range("D20").Value = 0.4986588806540624
$ws.Range("E20").Value = 0.1457289677622366
$ws.Range("G20").Value = 3.235958130036266
$ws.Range("H20").Value = 2.476790571514215
$ws.Range("J20").Value = 0.04939245777696044
$ws.Range("K20").Value = 1.243348297651295
$ws.Range("L20").Value = 0.4323064526603844
$ws.Range("M20").Value = 0.4747036361317711
$ws.Range("N20").Value = 4.214311391330142
$ws.Range("B21").Value = 1.930119746769492
$ws.Range("C21").Value = 0.04847081704757272
$ws.Range("D21").Value = 0.5018869638032584
$ws.Range("E21").Value = 0.1454921230329767
$ws.Range("G21").Value = 3.250031605745136
$ws.Range("H21").Value = 2.477009644577691
$ws.Range("J21").Value = 0.04965342112025439
$ws.Range("K21").Value = 1.29146391896748
$ws.Range("L21").Value = 0.4351954741713229
$ws.Range("M21").Value = 0.4833451950242775
$ws.Range("N21").Value = 4.19658188906179
$ws.Range("B22").Value = 1.963150968704326
$ws.Range("C22").Value = 0.05097726973670547
$ws.Range("D22").Value = 0.5041600490828699
$ws.Range("E22").Value = 0.1453690105602305
$ws.Range("G22").Value = 3.260282526780713
$ws.Range("H22").Value = 2.477825249510602
$ws.Range("J22").Value = 0.0498207445943919
$ws.Range("K22").Value = 1.323465305427248
$ws.Range("L22").Value = 0.437226068162019
$ws.Range("M22").Value = 0.4891801247275396
$ws.Range("N22").Value = 4.185785335855556
$ws.Range("B23").Value = 1.945449367595245
$ws.Range("C23").Value = 0.04964071760099387
$ws.Range("D23").Value = 0.5029319389169586
$ws.Range("E23").Value = 0.1454318255307587
$ws.Range("G23").Value = 3.254715062348737
$ws.Range("H23").Value = 2.477328463176804
$ws.Range("J23").Value = 0.04973173527744379
$ws.Range("K23").Value = 1.306334845735734
$ws.Range("L23").Value = 0.4361292907050114
$ws.Range("M23").Value = 0.4860488162058658
$ws.Range("N23").Value = 4.191475303414606
$ws.Range("B24").Value = 1.879995873607243
$ws.Range("C24").Value = 0.04455246634448429
$ws.Range("D24").Value = 0.4986088824072965
$ws.Range("E24").Value = 0.1457333750625889
$ws.Range("G24").Value = 3.235745885505651
$ws.Range("H24").Value = 2.476798352888551
$ws.Range("J24").Value = 0.04938813763146754
$ws.Range("K24").Value = 1.24257125863295
$ws.Range("L24").Value = 0.4322616431508237
$ws.Range("M24").Value = 0.4745655607746713
$ws.Range("N24").Value = 4.214614518125543
$ws.Range("B25").Value = 1.812800154358541
$ws.Range("C25").Value = 0.0390068512662225
$ws.Range("D25").Value = 0.4946495295473454
$ws.Range("E25").Value = 0.1461952829916697
$ws.Range("G25").Value = 3.219806142817021
$ws.Range("H25").Value = 2.479118065668985
$ws.Range("J25").Value = 0.04900350864449621
$ws.Range("K25").Value = 1.176181720944101
$ws.Range("L25").Value = 0.4287036813624923
$ws.Range("M25").Value = 0.4629855504899751
$ws.Range("N25").Value = 4.242968152409446
